# Daily attendance processing - 2026-01-09 16:40:12
# Toggles the order of names in the "Recorded By" column (G) wherever both
# "System" and "dnasr281@gmail.com" are recorded together for a session,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" and
# vice versa. Rows recorded by a single party are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nameA = "System"
$nameB = "dnasr281@gmail.com"
$comboAB = "$nameA, $nameB"
$comboBA = "$nameB, $nameA"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $comboAB) {
        $cell.Value2 = $comboBA
    }
    elseif ($val -eq $comboBA) {
        $cell.Value2 = $comboAB
    }
}
